$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one
# day (2026-02-08 -> 2026-02-09, i.e. 46061 -> 46062) for every data row
# (rows 2 through 374).
$ws.Range("C2:C374").Value = 46062
